$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 1.5
$ws.Cells.Item(2, 7).Value = 1.59
$ws.Cells.Item(2, 8).Value = 8.199999999999999
$ws.Cells.Item(2, 10).Value = 3.95
$ws.Cells.Item(2, 11).Value = 4.4
$ws.Cells.Item(2, 13).Value = 1.09
$ws.Cells.Item(2, 14).Value = 2.98
$ws.Cells.Item(2, 15).Value = 1.42
$ws.Cells.Item(2, 16).Value = 1.65
$ws.Cells.Item(2, 17).Value = 2.26
$ws.Cells.Item(2, 18).Value = 1.24
$ws.Cells.Item(2, 19).Value = 4.3
$ws.Cells.Item(2, 20).Value = 2.34
$ws.Cells.Item(2, 21).Value = 1.6
$ws.Cells.Item(2, 23).Value = 2.68
$ws.Cells.Item(2, 24).Value = 1000
$ws.Cells.Item(2, 28).Value = 6.4
$ws.Cells.Item(2, 31).Value = 240
$ws.Cells.Item(2, 33).Value = 1000
$ws.Cells.Item(2, 35).Value = 220
$ws.Cells.Item(2, 39).Value = 330
$ws.Cells.Item(2, 40).Value = 14
$ws.Cells.Item(2, 41).Value = 1000

# Row 3
$ws.Cells.Item(3, 7).Value = 4.4
$ws.Cells.Item(3, 8).Value = 1.91
$ws.Cells.Item(3, 9).Value = 2.04
$ws.Cells.Item(3, 10).Value = 3.8
$ws.Cells.Item(3, 11).Value = 4.3
$ws.Cells.Item(3, 18).Value = 1.62
$ws.Cells.Item(3, 19).Value = 2.32
$ws.Cells.Item(3, 22).Value = 1.96
$ws.Cells.Item(3, 24).Value = 29
$ws.Cells.Item(3, 25).Value = 14.5
$ws.Cells.Item(3, 26).Value = 16.5
$ws.Cells.Item(3, 27).Value = 29
$ws.Cells.Item(3, 28).Value = 28
$ws.Cells.Item(3, 29).Value = 12
$ws.Cells.Item(3, 30).Value = 13
$ws.Cells.Item(3, 31).Value = 22
$ws.Cells.Item(3, 32).Value = 40
$ws.Cells.Item(3, 33).Value = 22
$ws.Cells.Item(3, 34).Value = 18
$ws.Cells.Item(3, 35).Value = 30
$ws.Cells.Item(3, 36).Value = 75
$ws.Cells.Item(3, 37).Value = 50
$ws.Cells.Item(3, 38).Value = 50
$ws.Cells.Item(3, 40).Value = 30
$ws.Cells.Item(3, 41).Value = 10.5

# Row 4
$ws.Cells.Item(4, 12).Value = 1.38
$ws.Cells.Item(4, 13).Value = 1.06
$ws.Cells.Item(4, 14).Value = 4.4
$ws.Cells.Item(4, 15).Value = 1.28
$ws.Cells.Item(4, 16).Value = 2.14
$ws.Cells.Item(4, 17).Value = 1.84
$ws.Cells.Item(4, 20).Value = 1.77
$ws.Cells.Item(4, 21).Value = 2.2
$ws.Cells.Item(4, 24).Value = 16
$ws.Cells.Item(4, 28).Value = 18.5
$ws.Cells.Item(4, 30).Value = 9.800000000000001
$ws.Cells.Item(4, 34).Value = 18
$ws.Cells.Item(4, 36).Value = 95
$ws.Cells.Item(4, 41).Value = 11.5

# Row 5
$ws.Cells.Item(5, 8).Value = 2.22
$ws.Cells.Item(5, 12).Value = 1.32
$ws.Cells.Item(5, 14).Value = 5.4
$ws.Cells.Item(5, 15).Value = 1.21
$ws.Cells.Item(5, 16).Value = 2.44
$ws.Cells.Item(5, 17).Value = 1.66
$ws.Cells.Item(5, 18).Value = 1.58
$ws.Cells.Item(5, 19).Value = 2.62
$ws.Cells.Item(5, 20).Value = 1.57
$ws.Cells.Item(5, 24).Value = 21
$ws.Cells.Item(5, 25).Value = 14
$ws.Cells.Item(5, 26).Value = 16
$ws.Cells.Item(5, 27).Value = 28
$ws.Cells.Item(5, 31).Value = 20
$ws.Cells.Item(5, 35).Value = 28
$ws.Cells.Item(5, 40).Value = 23

# Row 6
$ws.Cells.Item(6, 10).Value = 3.5
$ws.Cells.Item(6, 14).Value = 3.95
$ws.Cells.Item(6, 15).Value = 1.26
$ws.Cells.Item(6, 17).Value = 1.76
$ws.Cells.Item(6, 18).Value = 1.4
$ws.Cells.Item(6, 19).Value = 2.96
$ws.Cells.Item(6, 20).Value = 1.67
$ws.Cells.Item(6, 21).Value = 2.2
$ws.Cells.Item(6, 24).Value = 19.5
$ws.Cells.Item(6, 25).Value = 17
$ws.Cells.Item(6, 27).Value = 65
$ws.Cells.Item(6, 28).Value = 13
$ws.Cells.Item(6, 29).Value = 9.800000000000001

# Row 8
$ws.Cells.Item(8, 6).Value = 1.95
$ws.Cells.Item(8, 7).Value = 2.96
$ws.Cells.Item(8, 9).Value = 3.65
$ws.Cells.Item(8, 12).Value = 1.28
$ws.Cells.Item(8, 14).Value = 1.32
$ws.Cells.Item(8, 16).Value = 1.32
$ws.Cells.Item(8, 18).Value = 1.32
$ws.Cells.Item(8, 19).Value = 2.6
$ws.Cells.Item(8, 22).Value = 1.37
$ws.Cells.Item(8, 23).Value = 1.51

# Row 9
$ws.Cells.Item(9, 13).Value = 1.15
$ws.Cells.Item(9, 14).Value = 2.38
$ws.Cells.Item(9, 15).Value = 1.64
$ws.Cells.Item(9, 16).Value = 1.45
$ws.Cells.Item(9, 19).Value = 6.2
$ws.Cells.Item(9, 20).Value = 2.28
$ws.Cells.Item(9, 21).Value = 1.68
$ws.Cells.Item(9, 24).Value = 7.6
$ws.Cells.Item(9, 25).Value = 10
$ws.Cells.Item(9, 27).Value = 130
$ws.Cells.Item(9, 28).Value = 6.8
$ws.Cells.Item(9, 29).Value = 8.6
$ws.Cells.Item(9, 30).Value = 970
$ws.Cells.Item(9, 31).Value = 80
$ws.Cells.Item(9, 32).Value = 16
$ws.Cells.Item(9, 33).Value = 15.5
$ws.Cells.Item(9, 34).Value = 970
$ws.Cells.Item(9, 35).Value = 140
$ws.Cells.Item(9, 36).Value = 970
$ws.Cells.Item(9, 37).Value = 46
$ws.Cells.Item(9, 38).Value = 75
$ws.Cells.Item(9, 39).Value = 300
$ws.Cells.Item(9, 40).Value = 970
$ws.Cells.Item(9, 41).Value = 160

# Row 10
$ws.Cells.Item(10, 7).Value = 3.65
$ws.Cells.Item(10, 9).Value = 3.25
$ws.Cells.Item(10, 10).Value = 2.78
$ws.Cells.Item(10, 14).Value = 2.5
$ws.Cells.Item(10, 15).Value = 1.53
$ws.Cells.Item(10, 16).Value = 1.5
$ws.Cells.Item(10, 17).Value = 2.42
$ws.Cells.Item(10, 21).Value = 1.56
$ws.Cells.Item(10, 22).Value = 1.44
$ws.Cells.Item(10, 24).Value = 11

# Row 12
$ws.Cells.Item(12, 7).Value = 2.24
$ws.Cells.Item(12, 8).Value = 4.7
$ws.Cells.Item(12, 9).Value = 5.8
$ws.Cells.Item(12, 12).Value = 1.01
$ws.Cells.Item(12, 13).Value = 1.15
$ws.Cells.Item(12, 14).Value = 2.32
$ws.Cells.Item(12, 15).Value = 1.65
$ws.Cells.Item(12, 18).Value = 1.14
$ws.Cells.Item(12, 19).Value = 6.4
$ws.Cells.Item(12, 20).Value = 2.32
$ws.Cells.Item(12, 21).Value = 1.6
$ws.Cells.Item(12, 22).Value = 1.21
$ws.Cells.Item(12, 23).Value = 1.81
$ws.Cells.Item(12, 24).Value = 7.6
$ws.Cells.Item(12, 25).Value = 970
$ws.Cells.Item(12, 26).Value = 970
$ws.Cells.Item(12, 27).Value = 190
$ws.Cells.Item(12, 28).Value = 6.4
$ws.Cells.Item(12, 29).Value = 7.6
$ws.Cells.Item(12, 30).Value = 28
$ws.Cells.Item(12, 31).Value = 130
$ws.Cells.Item(12, 32).Value = 970
$ws.Cells.Item(12, 33).Value = 970
$ws.Cells.Item(12, 34).Value = 970
$ws.Cells.Item(12, 35).Value = 170
$ws.Cells.Item(12, 36).Value = 970
$ws.Cells.Item(12, 37).Value = 970
$ws.Cells.Item(12, 38).Value = 80
$ws.Cells.Item(12, 39).Value = 330
$ws.Cells.Item(12, 40).Value = 970
$ws.Cells.Item(12, 41).Value = 240

# Row 13
$ws.Cells.Item(13, 7).Value = 3.1
$ws.Cells.Item(13, 8).Value = 3.1
$ws.Cells.Item(13, 9).Value = 3.4
$ws.Cells.Item(13, 11).Value = 2.9
$ws.Cells.Item(13, 12).Value = 1.01
$ws.Cells.Item(13, 13).Value = 1.01
$ws.Cells.Item(13, 14).Value = 2
$ws.Cells.Item(13, 15).Value = 1.79
$ws.Cells.Item(13, 18).Value = 1.1
$ws.Cells.Item(13, 19).Value = 8.4
$ws.Cells.Item(13, 20).Value = 2.28
$ws.Cells.Item(13, 21).Value = 1.49
$ws.Cells.Item(13, 22).Value = 1.42
$ws.Cells.Item(13, 23).Value = 1.47
$ws.Cells.Item(13, 24).Value = 7
$ws.Cells.Item(13, 25).Value = 8.4
$ws.Cells.Item(13, 26).Value = 970
$ws.Cells.Item(13, 27).Value = 90
$ws.Cells.Item(13, 28).Value = 970
$ws.Cells.Item(13, 29).Value = 8
$ws.Cells.Item(13, 30).Value = 970
$ws.Cells.Item(13, 31).Value = 80
$ws.Cells.Item(13, 32).Value = 24
$ws.Cells.Item(13, 33).Value = 970
$ws.Cells.Item(13, 34).Value = 970
$ws.Cells.Item(13, 35).Value = 1000
$ws.Cells.Item(13, 36).Value = 80
$ws.Cells.Item(13, 37).Value = 80
$ws.Cells.Item(13, 38).Value = 1000
$ws.Cells.Item(13, 39).Value = 1000
$ws.Cells.Item(13, 40).Value = 1000
$ws.Cells.Item(13, 41).Value = 1000

# Row 14
$ws.Cells.Item(14, 6).Value = 1.92
$ws.Cells.Item(14, 7).Value = 2.12
$ws.Cells.Item(14, 9).Value = 6.2
$ws.Cells.Item(14, 11).Value = 3.6
$ws.Cells.Item(14, 12).Value = 1.01
$ws.Cells.Item(14, 13).Value = 1.08
$ws.Cells.Item(14, 14).Value = 2.9
$ws.Cells.Item(14, 15).Value = 1.43
$ws.Cells.Item(14, 16).Value = 1.6
$ws.Cells.Item(14, 17).Value = 2.16
$ws.Cells.Item(14, 18).Value = 1.2
$ws.Cells.Item(14, 19).Value = 3.85
$ws.Cells.Item(14, 20).Value = 1.01
$ws.Cells.Item(14, 21).Value = 1.76
$ws.Cells.Item(14, 22).Value = 1.24
$ws.Cells.Item(14, 23).Value = 1.9
$ws.Cells.Item(14, 24).Value = 970
$ws.Cells.Item(14, 25).Value = 19.5
$ws.Cells.Item(14, 26).Value = 50
$ws.Cells.Item(14, 27).Value = 1000
$ws.Cells.Item(14, 28).Value = 10
$ws.Cells.Item(14, 29).Value = 10
$ws.Cells.Item(14, 30).Value = 28
$ws.Cells.Item(14, 31).Value = 100
$ws.Cells.Item(14, 32).Value = 16
$ws.Cells.Item(14, 33).Value = 15
$ws.Cells.Item(14, 34).Value = 30
$ws.Cells.Item(14, 35).Value = 1000
$ws.Cells.Item(14, 36).Value = 36
$ws.Cells.Item(14, 37).Value = 36
$ws.Cells.Item(14, 38).Value = 70
$ws.Cells.Item(14, 39).Value = 1000
$ws.Cells.Item(14, 40).Value = 1000
$ws.Cells.Item(14, 41).Value = 1000
